$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.295968333333334
$ws.Range("H2").Value = 18.887905
$ws.Range("I2").Value = 0.5052862712055841
$ws.Range("J2").Value = 0.5052862712055841
$ws.Range("M2").Value = 18.43631966666667
$ws.Range("N2").Value = 55.308959
$ws.Range("O2").Value = 0.6034704469962782
$ws.Range("P2").Value = 0.6034704469962781
$ws.Range("Q2").Value = 116.0744848045439
$ws.Range("R2").Value = 1044.670363240895
$ws.Range("S2").Value = 0.3049253319455165
$ws.Range("T2").Value = 0.3049253319455164
$ws.Range("G3").Value = 6.295968333333334
$ws.Range("H3").Value = 18.887905
$ws.Range("I3").Value = 0.5052862712055841
$ws.Range("J3").Value = 0.5052862712055841
$ws.Range("O3").Value = 0.1750419652256785
$ws.Range("P3").Value = 0.1750419652256784
$ws.Range("Q3").Value = 33.66843568541945
$ws.Range("R3").Value = 303.015921168775
$ws.Range("S3").Value = 0.0884463019133806
$ws.Range("T3").Value = 0.08844630191338058
$ws.Range("G4").Value = 6.295968333333334
$ws.Range("H4").Value = 18.887905
$ws.Range("I4").Value = 0.5052862712055841
$ws.Range("J4").Value = 0.5052862712055841
$ws.Range("M4").Value = 6.766555
$ws.Range("N4").Value = 20.299665
$ws.Range("O4").Value = 0.2214875877780434
$ws.Range("P4").Value = 0.2214875877780434
$ws.Range("Q4").Value = 42.60201600575834
$ws.Range("R4").Value = 383.418144051825
$ws.Range("S4").Value = 0.1119146373466871
$ws.Range("T4").Value = 0.111914637346687
$ws.Range("I5").Value = 0.2025983155648483
$ws.Range("J5").Value = 0.2025983155648483
$ws.Range("M5").Value = 18.43631966666667
$ws.Range("N5").Value = 55.308959
$ws.Range("O5").Value = 0.6034704469962782
$ws.Range("P5").Value = 0.6034704469962781
$ws.Range("Q5").Value = 46.54093420220811
$ws.Range("R5").Value = 418.868407819873
$ws.Range("S5").Value = 0.122262096054612
$ws.Range("T5").Value = 0.122262096054612
$ws.Range("I6").Value = 0.2025983155648483
$ws.Range("J6").Value = 0.2025983155648483
$ws.Range("O6").Value = 0.1750419652256785
$ws.Range("P6").Value = 0.1750419652256784
$ws.Range("S6").Value = 0.03546320730788321
$ws.Range("T6").Value = 0.03546320730788321
$ws.Range("I7").Value = 0.2025983155648483
$ws.Range("J7").Value = 0.2025983155648483
$ws.Range("M7").Value = 6.766555
$ws.Range("N7").Value = 20.299665
$ws.Range("O7").Value = 0.2214875877780434
$ws.Range("P7").Value = 0.2214875877780434
$ws.Range("Q7").Value = 17.08159745136167
$ws.Range("R7").Value = 153.734377062255
$ws.Range("S7").Value = 0.04487301220235308
$ws.Range("T7").Value = 0.04487301220235308
$ws.Range("G8").Value = 3.639816666666666
$ws.Range("H8").Value = 10.91945
$ws.Range("I8").Value = 0.2921154132295675
$ws.Range("J8").Value = 0.2921154132295676
$ws.Range("M8").Value = 18.43631966666667
$ws.Range("N8").Value = 55.308959
$ws.Range("O8").Value = 0.6034704469962782
$ws.Range("P8").Value = 0.6034704469962781
$ws.Range("Q8").Value = 67.10482359472778
$ws.Range("R8").Value = 603.94341235255
$ws.Range("S8").Value = 0.1762830189961496
$ws.Range("T8").Value = 0.1762830189961496
$ws.Range("G9").Value = 3.639816666666666
$ws.Range("H9").Value = 10.91945
$ws.Range("I9").Value = 0.2921154132295675
$ws.Range("J9").Value = 0.2921154132295676
$ws.Range("O9").Value = 0.1750419652256785
$ws.Range("P9").Value = 0.1750419652256784
$ws.Range("Q9").Value = 19.46435033663889
$ws.Range("R9").Value = 175.17915302975
$ws.Range("S9").Value = 0.05113245600441466
$ws.Range("T9").Value = 0.05113245600441466
$ws.Range("G10").Value = 3.639816666666666
$ws.Range("H10").Value = 10.91945
$ws.Range("I10").Value = 0.2921154132295675
$ws.Range("J10").Value = 0.2921154132295676
$ws.Range("M10").Value = 6.766555
$ws.Range("N10").Value = 20.299665
$ws.Range("O10").Value = 0.2214875877780434
$ws.Range("P10").Value = 0.2214875877780434
$ws.Range("Q10").Value = 24.62901966491667
$ws.Range("R10").Value = 221.66117698425
$ws.Range("S10").Value = 0.06469993822900326
$ws.Range("T10").Value = 0.06469993822900326
